# TypeScript.pptx edit:
#  - The old last slide ("Questions?") is kept, but moved to become the new
#    last (10th) slide.
#  - The slide that used to be "Questions?" (slide 9) is turned into a new
#    "What next?" slide with a title + content placeholder containing two
#    hyperlinked references.

$p = $ppt.ActivePresentation

# --- Step 1: duplicate the current "Questions?" slide (slide 9). -----------
# Duplicate() inserts the copy immediately after the original, so it becomes
# the new slide 10 - exactly the "Questions?" slide we want to keep at the
# end of the deck.
$questionsSlide = $p.Slides.Item(9)
$questionsDup = $questionsSlide.Duplicate()

# --- Step 2: remove the original slide 9 and replace it with a fresh one ---
# using the "Title and Content" layout (title + body placeholder), matching
# the layout used by the other content slides in this deck.
$p.Slides.Item(9).Delete()
$whatNext = $p.Slides.Add(9, 16)

# --- Step 3: title text -----------------------------------------------------
$whatNext.Shapes.Item(1).TextFrame.TextRange.Text = "What next?"

# --- Step 4: body content with hyperlinks -----------------------------------
$docsUrl = "https://www.typescriptlang.org/docs/home.html"
$examplesUrl = "https://github.com/Mishurin/frontend-examples"

$body = $whatNext.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Read official documentation: " + $docsUrl + "`r`rExamples for this presentation are here: " + $examplesUrl + "`r"

$link1 = $body.Find($docsUrl)
$link1.ActionSettings.Item(1).Hyperlink.Address = $docsUrl

$link2 = $body.Find($examplesUrl)
$link2.ActionSettings.Item(1).Hyperlink.Address = $examplesUrl

# Last (4th) paragraph is blank with no bullet, same as in the authored deck.
$lastPara = $body.Paragraphs(4, 1)
$lastPara.ParagraphFormat.Bullet.Visible = 0
